$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stray row that held "5963230 - Leandro Gonçalves de Aguiar" on
# its own (old row 13); everything below shifts up by one row.
$ws.Rows(13).Delete()

# Objetivos: now shows the docente responsável instead of the old long text.
$ws.Range("B10").Value = "5963230 - Leandro Gonçalves de Aguiar"
$ws.Range("C10").Value = "5963230 - Leandro Gonçalves de Aguiar"

# Programa resumido: now just "Semestral".
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# Programa: now just the activation date "01/01/2012" (same text already
# used by the Ativação row, B8:C8) -- copy it across so it lands as a text
# shared-string instead of being auto-parsed into a date serial number.
$ws.Range("B8:C8").Copy()
$ws.Range("B15").PasteSpecial(-4163)
$excel.CutCopyMode = $false

# Método: now shows the docente responsável.
$ws.Range("B18").Value = "5963230 - Leandro Gonçalves de Aguiar"
$ws.Range("C18").Value = "5963230 - Leandro Gonçalves de Aguiar"

# Critério: now shows the old "Método" description text.
$ws.Range("B19").Value = "O desenvolvimento da disciplina será baseado em leituras, aula expositiva, discussão e resolução de estudos de caso e resolução de exercícios."
$ws.Range("C19").Value = "O desenvolvimento da disciplina será baseado em leituras, aula expositiva, discussão e resolução de estudos de caso e resolução de exercícios."

# Norma de recuperação: now shows the old "Critério" text.
$ws.Range("B20").Value = "Provas e trabalhos."
$ws.Range("C20").Value = "Provas e trabalhos."

# Bibliografia: now shows the old "Norma de recuperação" text.
$ws.Range("B21").Value = "Prova única com nota maior ou igual a 5,0 (cinco)."
$ws.Range("C21").Value = "Prova única com nota maior ou igual a 5,0 (cinco)."
